$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($row, $col, $value) {
    $c = $ws.Cells.Item($row, $col)
    $c.NumberFormat = "@"
    $c.Value = $value
    $c.ClearFormats()
}

Set-TextValue 2 4 "26.774.42"
Set-TextValue 2 5 "  +1.31%  "
Set-TextValue 3 4 "1.725.80"
Set-TextValue 3 5 "  +0.06%  "
Set-TextValue 4 5 "  -0.23%  "
Set-TextValue 5 4 "240.82"
Set-TextValue 5 5 "  -1.03%  "
Set-TextValue 6 4 "0.9979"
Set-TextValue 6 5 "  -0.21%  "
Set-TextValue 7 4 "0.4857"
Set-TextValue 7 5 "  -0.82%  "
Set-TextValue 8 4 "0.2582"
Set-TextValue 8 5 "  -0.84%  "
Set-TextValue 9 4 "0.06192"
Set-TextValue 9 5 "  -0.15%  "
Set-TextValue 10 4 "1.727.58"
Set-TextValue 10 5 "  +0.26%  "
Set-TextValue 11 4 "15.93"
Set-TextValue 11 5 "  +2.91%  "
Set-TextValue 12 4 "0.06882"
Set-TextValue 12 5 "  -1.52%  "
Set-TextValue 13 4 "0.6078"
Set-TextValue 13 5 "  +1.44%  "
Set-TextValue 14 4 "4.472"
Set-TextValue 14 5 "  -1.44%  "
Set-TextValue 15 4 "76.96"
Set-TextValue 15 5 "  -0.57%  "
Set-TextValue 16 5 "  -0.19%  "
Set-TextValue 17 4 "26.576.64"
Set-TextValue 17 5 "  +0.58%  "
Set-TextValue 19 4 "0.000007147"
Set-TextValue 19 5 "  -1.19%  "
Set-TextValue 20 4 "11.41"
Set-TextValue 20 5 "  +0.57%  "
Set-TextValue 21 4 "1.950.30"
Set-TextValue 21 5 "  +0.32%  "
Set-TextValue 22 4 "4.428"
Set-TextValue 22 5 "  -0.77%  "
Set-TextValue 23 4 "8.575"
Set-TextValue 23 5 "  +0.03%  "
Set-TextValue 24 4 "5.089"
Set-TextValue 24 5 "  -1.17%  "
Set-TextValue 25 4 "137.47"
Set-TextValue 25 5 "  -0.28%  "
Set-TextValue 26 4 "15.23"
Set-TextValue 26 5 "  -0.34%  "
Set-TextValue 27 4 "1.775"
Set-TextValue 27 5 "  +2.75%  "
Set-TextValue 28 4 "105.94"
Set-TextValue 28 5 "  -0.91%  "
Set-TextValue 29 5 "  -1.93%  "
Set-TextValue 30 4 "3.994"
Set-TextValue 30 5 "  +1.08%  "
Set-TextValue 31 4 "0.07927"
Set-TextValue 32 4 "3.683"
Set-TextValue 33 5 "  -0.60%  "
Set-TextValue 34 2 "Frax"
Set-TextValue 34 3 "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
Set-TextValue 34 4 "0.9970"
Set-TextValue 34 5 "  -0.22%  "
Set-TextValue 35 2 "HuobiToken"
Set-TextValue 35 3 "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
Set-TextValue 35 4 "2.596"
Set-TextValue 35 5 "  -0.25%  "
Set-TextValue 36 2 "ARBITRUM"
Set-TextValue 36 3 "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
Set-TextValue 36 4 "1.005"
Set-TextValue 36 5 "  +0.13%  "
Set-TextValue 37 2 "ImmutableX"
Set-TextValue 37 3 "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
Set-TextValue 37 4 "0.6206"
Set-TextValue 37 5 "  -0.93%  "
Set-TextValue 38 2 "TrustWalletToken"
Set-TextValue 38 3 "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
Set-TextValue 38 4 "0.9226"
Set-TextValue 38 5 "  -1.87%  "
Set-TextValue 39 2 "RenderToken"
Set-TextValue 39 3 "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
Set-TextValue 39 4 "2.024"
Set-TextValue 39 5 "  +3.96%  "
Set-TextValue 40 2 "MXToken"
Set-TextValue 40 3 "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
Set-TextValue 40 4 "2.445"
Set-TextValue 40 5 "  +2.41%  "
Set-TextValue 41 2 "PaxDollar"
Set-TextValue 41 3 "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
Set-TextValue 41 4 "0.9972"
Set-TextValue 41 5 "  -0.26%  "
Set-TextValue 42 2 "VeChain"
Set-TextValue 42 3 "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
Set-TextValue 42 4 "0.01495"
Set-TextValue 42 5 "  +0.80%  "
Set-TextValue 43 2 "FraxShare"
Set-TextValue 43 3 "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
Set-TextValue 43 4 "5.617"
Set-TextValue 43 5 "  +5.94%  "
Set-TextValue 44 2 "Quant"
Set-TextValue 44 3 "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
Set-TextValue 44 4 "99.82"
Set-TextValue 44 5 "  +0.33%  "
Set-TextValue 45 2 "TheSandbox"
Set-TextValue 45 3 "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
Set-TextValue 45 4 "0.3839"
Set-TextValue 45 5 "  -0.37%  "
Set-TextValue 46 2 "Aptos"
Set-TextValue 46 3 "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
Set-TextValue 46 4 "6.838"
Set-TextValue 46 5 "  +0.05%  "
Set-TextValue 47 2 "Algorand"
Set-TextValue 47 3 "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
Set-TextValue 47 4 "0.1155"
Set-TextValue 47 5 "  -1.30%  "
Set-TextValue 48 2 "Cronos"
Set-TextValue 48 3 "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
Set-TextValue 48 4 "0.05377"
Set-TextValue 48 5 "  +0.24%  "
Set-TextValue 49 2 "EnergySwap"
Set-TextValue 49 3 "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
Set-TextValue 49 4 "7.859"
Set-TextValue 49 5 "  +1.40%  "
Set-TextValue 50 2 "Elrond"
Set-TextValue 50 3 "https://coinranking.com/coin/omwkOTglq+elrond-egld"
Set-TextValue 50 4 "30.09"
Set-TextValue 50 5 "  -0.44%  "
Set-TextValue 51 2 "NEARProtocol"
Set-TextValue 51 3 "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
Set-TextValue 51 4 "1.233"
Set-TextValue 51 5 "  -0.24%  "
